$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5040.6313
$ws.Range("I40").Value = 4505.6665
$ws.Range("J40").Value = 5522.1
$ws.Range("K40").Value = 4505.6665
$ws.Range("L40").Value = 5522.1
$ws.Range("M40").Value = -4330.6665
$ws.Range("N40").Value = -5872.1
$ws.Range("H92").Value = 126
$ws.Range("I92").Value = 39.5
$ws.Range("K92").Value = 39.5
$ws.Range("M92").Value = 1208.5
$ws.Range("H112").Value = 3498.6
$ws.Range("J112").Value = 3498.6
$ws.Range("L112").Value = 10495.8
$ws.Range("N112").Value = -12711.8
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = ""
$ws.Range("H127").Value = 4984.5
$ws.Range("I127").Value = 4979.3335
$ws.Range("K127").Value = 14938.0005
$ws.Range("M127").Value = -9978.000499999998
$ws.Range("H132").Value = 9035.200000000001
$ws.Range("I132").Value = 8252.036
$ws.Range("K132").Value = 24756.108
$ws.Range("M132").Value = -22226.108
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1440.7142
$ws.Range("I74").Value = 1222.5
$ws.Range("K74").Value = 1222.5
$ws.Range("M74").Value = -348.5
$ws.Range("H77").Value = 1440.7142
$ws.Range("I77").Value = 1222.5
$ws.Range("K77").Value = 6112.5
$ws.Range("M77").Value = -1744.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 421
$ws.Range("I12").Value = 421
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 421
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -253
$ws.Range("N12").Value = ""
$ws.Range("H62").Value = 97590.5
$ws.Range("J62").Value = 97590.5
$ws.Range("L62").Value = 97590.5
$ws.Range("N62").Value = -98962.5
$ws.Range("H65").Value = 97590.5
$ws.Range("J65").Value = 97590.5
$ws.Range("L65").Value = 292771.5
$ws.Range("N65").Value = -299635.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1075
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 1133.3334
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 3400.0002
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -5022.0002
$ws.Range("H71").Value = 1075
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 1133.3334
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 10200.0006
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -18312.0006
$ws.Range("H113").Value = 2156.96
$ws.Range("I113").Value = 1110
$ws.Range("J113").Value = 2487.5789
$ws.Range("K113").Value = 3330
$ws.Range("L113").Value = 7462.736699999999
$ws.Range("M113").Value = -1160
$ws.Range("N113").Value = -11802.7367
$ws.Range("H132").Value = 3044.95
$ws.Range("I132").Value = 2699.8572
$ws.Range("J132").Value = 3230.7693
$ws.Range("K132").Value = 24298.7148
$ws.Range("L132").Value = 29076.9237
$ws.Range("M132").Value = -21768.7148
$ws.Range("N132").Value = -34136.9237
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8000
$ws.Range("I5").Value = 7000
$ws.Range("J5").Value = 9000
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = -6888
$ws.Range("N5").Value = -9224
$ws.Range("H13").Value = 401
$ws.Range("I13").Value = 275
$ws.Range("J13").Value = 485
$ws.Range("K13").Value = 275
$ws.Range("L13").Value = 485
$ws.Range("M13").Value = -136
$ws.Range("N13").Value = -763
$ws.Range("H113").Value = 7558.7646
$ws.Range("I113").Value = 6066.5557
$ws.Range("J113").Value = 9237.5
$ws.Range("K113").Value = 6066.5557
$ws.Range("L113").Value = 9237.5
$ws.Range("M113").Value = -3896.5557
$ws.Range("N113").Value = -13577.5
$ws.Range("H132").Value = 103818.55
$ws.Range("I132").Value = 139250.88
$ws.Range("K132").Value = 417752.64
$ws.Range("M132").Value = -415222.64
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4000
$ws.Range("I2").Value = 4000
$ws.Range("K2").Value = 4000
$ws.Range("M2").Value = -3888
$ws.Range("H40").Value = 9453.546
$ws.Range("I40").Value = 9141.429
$ws.Range("K40").Value = 9141.429
$ws.Range("M40").Value = -9005.429
$ws.Range("N40").Value = -10271.75
$ws.Range("H61").Value = 4156.143
$ws.Range("I61").Value = 2955
$ws.Range("J61").Value = 7999.8
$ws.Range("K61").Value = 2955
$ws.Range("L61").Value = 7999.8
$ws.Range("M61").Value = -2753
$ws.Range("N61").Value = -8403.799999999999
$ws.Range("H113").Value = 4156.143
$ws.Range("I113").Value = 2955
$ws.Range("J113").Value = 7999.8
$ws.Range("K113").Value = 2955
$ws.Range("L113").Value = 7999.8
$ws.Range("M113").Value = -785
$ws.Range("N113").Value = -12339.8
$ws.Range("H132").Value = 6998.2856
$ws.Range("I132").Value = 4998.25
$ws.Range("J132").Value = 9665
$ws.Range("K132").Value = 14994.75
$ws.Range("L132").Value = 28995
$ws.Range("M132").Value = -12464.75
$ws.Range("N132").Value = -34055
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15500
$ws.Range("H36").Value = 15000
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15500
$ws.Range("H48").Value = 55000
$ws.Range("J48").Value = 55000
$ws.Range("L48").Value = 55000
$ws.Range("N48").Value = -56138
$ws.Range("H63").Value = 23374.834
$ws.Range("J63").Value = 25049.8
$ws.Range("L63").Value = 25049.8
$ws.Range("N63").Value = -26297.8
$ws.Range("H66").Value = 23374.834
$ws.Range("J66").Value = 25049.8
$ws.Range("L66").Value = 75149.39999999999
$ws.Range("N66").Value = -81389.39999999999
$ws.Range("H81").Value = 1933.5
$ws.Range("I81").Value = 1933.5
$ws.Range("K81").Value = 3867
$ws.Range("M81").Value = -2806
$ws.Range("H84").Value = 1933.5
$ws.Range("I84").Value = 1933.5
$ws.Range("K84").Value = 19335
$ws.Range("M84").Value = -14031
$ws.Range("H100").Value = 537
$ws.Range("I100").Value = 365.13333
$ws.Range("K100").Value = 730.26666
$ws.Range("M100").Value = -189.26666
$ws.Range("H113").Value = 549.3684
$ws.Range("I113").Value = 435.72726
$ws.Range("K113").Value = 1307.18178
$ws.Range("M113").Value = 862.8182200000001
$ws.Range("H132").Value = 2924.239
$ws.Range("I132").Value = 2673.4324
$ws.Range("J132").Value = 3955.3333
$ws.Range("K132").Value = 8020.297200000001
$ws.Range("L132").Value = 11865.9999
$ws.Range("M132").Value = -5490.297200000001
$ws.Range("N132").Value = -16925.9999
$ws.Range("H136").Value = 3853.311
$ws.Range("I136").Value = 2483.0293
$ws.Range("K136").Value = 7449.0879
$ws.Range("M136").Value = -4899.0879
